$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week")

$ws.Range("I19").Select()

$ws.Range("C19").Value = 7
$ws.Range("E19").Value = 4
$ws.Range("G19").Value = 6
$ws.Range("I19").Value = 6

$ws.Range("C18").Clear()
$ws.Range("E18").Clear()
$ws.Range("G18").Clear()
$ws.Range("I18").Clear()

$ws.Range("C21:I21").Formula = "=SUM(C15:C20)"

$wb.Save()
